# Remove `pax` from databases
# Rename header labels in the USE_TYPE_PROPERTIES workbook so the
# occupant-normalised units drop the redundant "pax" suffix.

$wb = $excel.ActiveWorkbook

# --- INTERNAL_LOADS sheet ---------------------------------------------
$ws1 = $wb.Worksheets.Item("INTERNAL_LOADS")

$ws1.Range("I1").Value = "Vww_ldp"
$ws1.Range("J1").Value = "Vw_ldp"
$ws1.Range("D1").Value = "X_ghp"
$ws1.Range("C1").Value = "Qs_Wp"
$ws1.Range("B1").Value = "Occ_m2p"

# --- INDOOR_COMFORT sheet ----------------------------------------------
$ws2 = $wb.Worksheets.Item("INDOOR_COMFORT")

$ws2.Range("F1").Value = "Ve_lsp"
